$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: reposition/resize the saved window to match the author's
# on-screen state. Harmless if unsupported by the host.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 11657
    $win.Top = 3429
    $win.Width = 17717
    $win.Height = 13225
} catch {
}

$rows = @(
    @(1, 919.115, 1030.576, 2550.4, 0),
    @(2, 941.478, 1802.631, 2550.4, -16.07),
    @(3, 932.207, 1620.287, 2550.4, -12),
    @(4, 928.866, 1542.607, 2550.4, -9.97),
    @(5, 924.67, 1412.415, 2550.4, -7.97),
    @(6, 920.349, 1218.783, 2550.4, -3.97),
    @(7, 919.523, 1124.577, 2550.4, -2),
    @(8, 919.105, 1030.13, 2550.4, 0),
    @(9, 919.971, 936.123, 2550.4, 2),
    @(10, 921, 841.029, 2550.4, 3.98),
    @(11, 922.817, 744.58, 2550.4, 5.98),
    @(12, 925.446, 647.851, 2550.4, 7.95),
    @(13, 928.813, 548.542, 2550.4, 9.98),
    @(14, 932.871, 450.189, 2550.4, 11.98),
    @(15, 937.574, 350.344, 2550.4, 13.98),
    @(16, 943.245, 247.597, 2550.4, 15.93),
    @(17, 901.075, 1030.555, 2650.1, 0),
    @(18, 922.543, 1801.894, 2650.1, -16.07),
    @(19, 912.638, 1599.435, 2650.1, -11.97),
    @(20, 909.189, 1502.908, 2650.1, -9.97),
    @(21, 906.088, 1408.439, 2650.1, -7.97),
    @(22, 903.999, 1312.886, 2650.1, -6),
    @(23, 902.209, 1218.653, 2650.1, -3.97),
    @(24, 901.514, 1124.627, 2650.1, -2),
    @(25, 901.082, 1030.42, 2650.1, 0),
    @(26, 901.838, 936.304, 2650.1, 2),
    @(27, 902.954, 841.531, 2650.1, 3.98),
    @(28, 904.666, 745.056, 2650.1, 5.98),
    @(29, 907.211, 648.743, 2650.1, 7.95),
    @(30, 909.73, 569.355, 2650.1, 9.98),
    @(31, 914.225, 450.832, 2650.1, 11.98),
    @(32, 918.984, 350.968, 2650.1, 13.98),
    @(33, 924.484, 248.432, 2650.1, 15.93),
    @(34, 884.5, 1030.753, 2749.8, 0),
    @(35, 905.445, 1801.417, 2749.8, -16.08),
    @(36, 900.359, 1698.944, 2749.8, -14),
    @(37, 895.807, 1599.66, 2749.8, -11.98),
    @(38, 892.462, 1504.194, 2749.8, -10),
    @(39, 889.41, 1408.297, 2749.8, -7.98),
    @(40, 885.824, 1219.739, 2749.8, -3.98),
    @(41, 884.527, 1031.932, 2749.8, 0),
    @(42, 885.009, 937.771, 2749.8, 1.97),
    @(43, 886.006, 841.864, 2749.8, 4),
    @(44, 887.994, 746.771, 2749.8, 5.97),
    @(45, 890.22, 650.168, 2749.8, 7.97),
    @(46, 893.517, 552.154, 2749.8, 9.97),
    @(47, 897.083, 453.975, 2749.8, 11.95),
    @(48, 901.861, 352.632, 2749.8, 13.97),
    @(49, 905.995, 270.607, 2749.8, 15.95),
    @(50, 869.115, 1030.66, 2849.5, 0),
    @(51, 889.445, 1801.016, 2849.5, -16.08),
    @(52, 884.452, 1697.618, 2849.5, -14),
    @(53, 880.075, 1599.131, 2849.5, -11.98),
    @(54, 876.501, 1503.406, 2849.5, -9.98),
    @(55, 873.885, 1408.65, 2849.5, -7.98),
    @(56, 871.904, 1314.378, 2849.5, -6),
    @(57, 870.073, 1219.333, 2849.5, -3.98),
    @(58, 869.41, 1125.542, 2849.5, -2),
    @(59, 869.5, 937.761, 2849.5, 2),
    @(60, 870.793, 842.851, 2849.5, 3.97),
    @(61, 872.076, 746.548, 2849.5, 5.97),
    @(62, 874.059, 669.743, 2849.5, 7.97),
    @(63, 877.908, 552.619, 2849.5, 9.97),
    @(64, 881.376, 454.048, 2849.5, 11.97),
    @(65, 885.892, 354.83, 2849.5, 13.97),
    @(66, 890.948, 251.958, 2849.5, 15.95),
    @(67, 854.679, 1030.881, 2949.2, 0),
    @(68, 874.441, 1799.556, 2949.2, -16.07),
    @(69, 869.47, 1696.839, 2949.2, -14),
    @(70, 865.46, 1598.211, 2949.2, -12),
    @(71, 861.926, 1501.767, 2949.2, -9.97),
    @(72, 859.241, 1406.746, 2949.2, -7.97),
    @(73, 857.143, 1312.456, 2949.2, -6),
    @(74, 855.627, 1218.286, 2949.2, -3.97),
    @(75, 854.9, 1125.131, 2949.2, -2),
    @(76, 854.618, 1031.007, 2949.2, 0),
    @(77, 855.056, 936.757, 2949.2, 2),
    @(78, 856.022, 842.164, 2949.2, 3.98),
    @(79, 857.963, 747.051, 2949.2, 5.98),
    @(80, 860.021, 650.049, 2949.2, 7.95),
    @(81, 863.04, 551.439, 2949.2, 9.98),
    @(82, 866.826, 452.906, 2949.2, 11.98),
    @(83, 870.992, 353.676, 2949.2, 13.98),
    @(84, 874.987, 272.222, 2949.2, 15.93)
)

foreach ($r in $rows) {
    $ws.Cells.Item($r[0], 1).Value = $r[1]
    $ws.Cells.Item($r[0], 2).Value = $r[2]
    $ws.Cells.Item($r[0], 3).Value = $r[3]
    $ws.Cells.Item($r[0], 4).Value = $r[4]
}

# Clear rows 85:90 (A:D) back to blank cells
$ws.Range("A85:D90").ClearContents()

# Remove rows 121:126 entirely (shrinks used range to E120)
$ws.Range("A121:E126").EntireRow.Delete()

# Update selection to match target view state
$ws.Range("D12").Select() | Out-Null
